$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74
$ws.Cells.Item($row, 1).Value = "2025-04-29 11:00:37"
$ws.Cells.Item($row, 2).Value = 248
